$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): move the two surviving headers that shift right,
#    rewrite the header that changed wording, and add the new headers that
#    describe the simplified diagnostic columns.
# ---------------------------------------------------------------------------

# Re-home the two untouched headers that slide from S1/T1 to W1/X1 *before*
# their old cells are overwritten, so the shared strings stay referenced.
$ws.Range("W1").Value = "Gestionnaires pas notifiés"
$ws.Range("X1").Value = "Identifiant import"

# Existing header whose wording changed.
$ws.Range("L1").Value = "Année du diagnostic"

# New diagnostic columns (also overwrites the now-relocated old BIo/Durable/
# Rouge/AOC/HVE headers in N1:R1).
$ws.Range("N1").Value = "Bio"
$ws.Range("O1").Value = "SIQO"
$ws.Range("P1").Value = "Environnement"
$ws.Range("Q1").Value = "Autre EGAlim"
$ws.Range("R1").Value = "Viandes vollailes total"
$ws.Range("S1").Value = "Viandes vollailles EGAlim"
$ws.Range("T1").Value = "Viandes vollailles provenance France"
$ws.Range("U1").Value = "Produits aquatiques total"
$ws.Range("V1").Value = "Produits aquatiques EGAlim"
$ws.Range("Y1").Value = "Statut de publication"

# ---------------------------------------------------------------------------
# 2. Data row (row 2): one simplified, automated-test diagnostic import.
# ---------------------------------------------------------------------------

$ws.Range("A2").Value = 73282932000074
$ws.Range("B2").Value = "Staff canteen"
$ws.Range("D2").Value = 54460
$ws.Range("F2").Value = 700
$ws.Range("J2").Value = "public"
$ws.Range("L2").Value = 2019
$ws.Range("M2").Value = 1000
$ws.Range("N2").Value = 500
$ws.Range("O2").Value = 100.1
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("X2").Value = "Automated test"

# ---------------------------------------------------------------------------
# 3. Header style: the header font grows from 10pt to 12pt bold. Apply the
#    size change to the original header cells, then stamp the same format
#    onto the newly added header cells so every header shares one style.
# ---------------------------------------------------------------------------

$ws.Range("A1:T1").Font.Size = 12
$ws.Range("A1").Copy() | Out-Null
$ws.Range("U1:Y1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Column widths for the final A:Y layout.
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 13.498697916666666
$ws.Columns.Item(2).ColumnWidth = 15.498697916666666
$ws.Columns.Item(3).ColumnWidth = 5.166666666666667
$ws.Columns.Item(4).ColumnWidth = 9.830729166666666
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws.Columns.Item(6).ColumnWidth = 10.166666666666666
$ws.Columns.Item(7).ColumnWidth = 7.330729166666667
$ws.Columns.Item(8).ColumnWidth = 16.830729166666668
$ws.Columns.Item(9).ColumnWidth = 13.998697916666666
$ws.Columns.Item(10).ColumnWidth = 17.498697916666668
$ws.Columns.Item(11).ColumnWidth = 18.666666666666668
$ws.Columns.Item(12).ColumnWidth = 17.330729166666668
$ws.Columns.Item(13).ColumnWidth = 4.498697916666667
$ws.Columns.Item(14).ColumnWidth = 3.3307291666666665
$ws.Columns.Item(15).ColumnWidth = 4.498697916666667
$ws.Columns.Item(16).ColumnWidth = 12.830729166666666
$ws.Columns.Item(17).ColumnWidth = 11.498697916666666
$ws.Columns.Item(18).ColumnWidth = 18.666666666666668
$ws.Columns.Item(19).ColumnWidth = 21.666666666666668
$ws.Columns.Item(20).ColumnWidth = 31.330729166666668
$ws.Columns.Item(21).ColumnWidth = 21.166666666666668
$ws.Columns.Item(22).ColumnWidth = 23.498697916666668
$ws.Columns.Item(23).ColumnWidth = 22.166666666666668
$ws.Columns.Item(24).ColumnWidth = 14.830729166666666
$ws.Columns.Item(25).ColumnWidth = 17.498697916666668

# ---------------------------------------------------------------------------
# 5. View state: leave the selection on the newly imported identifier cell.
# ---------------------------------------------------------------------------

$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("X2").Select() | Out-Null
